# Processes crbeCountries manytomany relation added
#
# - Column P (16) header text shortened: "Срок дейст. патента" -> "Срок дейст."
# - Column P width narrowed (no longer needs to fit "... патента")
# - Active selection moved to the (now narrower) P1 header cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the column P header text.
$ws.Range("P1").Value = "Срок дейст."

# Narrow column P to match the shorter header text.
$ws.Columns.Item(16).ColumnWidth = 11

# Reflect the new selection left behind after editing that header cell.
$ws.Range("P1").Select()
